$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "L1cam"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.603177
$ws.Range("H2").Value = 19.809531
$ws.Range("I2").Value = 0.5135477412645301
$ws.Range("J2").Value = 0.5135477412645302
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.885873333333334
$ws.Range("N2").Value = 8.657620000000001
$ws.Range("O2").Value = 0.3070415651026022
$ws.Range("P2").Value = 0.3070415651026022
$ws.Range("Q2").Value = 19.05593241958
$ws.Range("R2").Value = 171.50339177622
$ws.Range("S2").Value = 0.1576805022327675
$ws.Range("T2").Value = 0.1576805022327676

# row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "L1cam"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.603177
$ws.Range("H3").Value = 19.809531
$ws.Range("I3").Value = 0.5135477412645301
$ws.Range("J3").Value = 0.5135477412645302
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.165953666666667
$ws.Range("N3").Value = 9.497861
$ws.Range("O3").Value = 0.3368406220840099
$ws.Range("P3").Value = 0.3368406220840099
$ws.Range("Q3").Value = 20.905352434799
$ws.Range("R3").Value = 188.148171913191
$ws.Range("S3").Value = 0.1729837406373825
$ws.Range("T3").Value = 0.1729837406373825

# row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "L1cam"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.603177
$ws.Range("H4").Value = 19.809531
$ws.Range("I4").Value = 0.5135477412645301
$ws.Range("J4").Value = 0.5135477412645302
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.327024333333334
$ws.Range("N4").Value = 9.981073
$ws.Range("O4").Value = 0.3539776838580724
$ws.Range("P4").Value = 0.3539776838580724
$ws.Range("Q4").Value = 21.968930556307
$ws.Range("R4").Value = 197.720375006763
$ws.Range("S4").Value = 0.181784440003363
$ws.Range("T4").Value = 0.1817844400033631

# row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "L1cam"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.603177
$ws.Range("H5").Value = 19.809531
$ws.Range("I5").Value = 0.5135477412645301
$ws.Range("J5").Value = 0.5135477412645302
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.020115
$ws.Range("N5").Value = 0.060345
$ws.Range("O5").Value = 0.002140128955315263
$ws.Range("P5").Value = 0.002140128955315263
$ws.Range("Q5").Value = 0.132822905355
$ws.Range("R5").Value = 1.195406148195
$ws.Range("S5").Value = 0.001099058391016972
$ws.Range("T5").Value = 0.001099058391016972

# row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "L1cam"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3785896666666667
$ws.Range("H6").Value = 1.135769
$ws.Range("I6").Value = 0.02944398858046029
$ws.Range("J6").Value = 0.0294439885804603
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.885873333333334
$ws.Range("N6").Value = 8.657620000000001
$ws.Range("O6").Value = 0.3070415651026022
$ws.Range("P6").Value = 0.3070415651026022
$ws.Range("Q6").Value = 1.092561823308889
$ws.Range("R6").Value = 9.833056409780001
$ws.Range("S6").Value = 0.009040528336607675
$ws.Range("T6").Value = 0.009040528336607677

# row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "L1cam"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3785896666666667
$ws.Range("H7").Value = 1.135769
$ws.Range("I7").Value = 0.02944398858046029
$ws.Range("J7").Value = 0.0294439885804603
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.165953666666667
$ws.Range("N7").Value = 9.497861
$ws.Range("O7").Value = 0.3368406220840099
$ws.Range("P7").Value = 0.3368406220840099
$ws.Range("Q7").Value = 1.198597343345444
$ws.Range("R7").Value = 10.787376090109
$ws.Range("S7").Value = 0.009917931430076729
$ws.Range("T7").Value = 0.009917931430076731

# row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "L1cam"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3785896666666667
$ws.Range("H8").Value = 1.135769
$ws.Range("I8").Value = 0.02944398858046029
$ws.Range("J8").Value = 0.0294439885804603
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.327024333333334
$ws.Range("N8").Value = 9.981073
$ws.Range("O8").Value = 0.3539776838580724
$ws.Range("P8").Value = 0.3539776838580724
$ws.Range("Q8").Value = 1.259577033348556
$ws.Range("R8").Value = 11.336193300137
$ws.Range("S8").Value = 0.01042251488125487
$ws.Range("T8").Value = 0.01042251488125487

# row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "L1cam"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3785896666666667
$ws.Range("H9").Value = 1.135769
$ws.Range("I9").Value = 0.02944398858046029
$ws.Range("J9").Value = 0.0294439885804603
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.020115
$ws.Range("N9").Value = 0.060345
$ws.Range("O9").Value = 0.002140128955315263
$ws.Range("P9").Value = 0.002140128955315263
$ws.Range("Q9").Value = 0.007615331145
$ws.Range("R9").Value = 0.06853798030500001
$ws.Range("S9").Value = 0.00006301393252101501
$ws.Range("T9").Value = 0.00006301393252101504

# row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "L1cam"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.441487333333333
$ws.Range("H10").Value = 10.324462
$ws.Range("I10").Value = 0.2676541983690312
$ws.Range("J10").Value = 0.2676541983690313
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.885873333333334
$ws.Range("N10").Value = 8.657620000000001
$ws.Range("O10").Value = 0.3070415651026022
$ws.Range("P10").Value = 0.3070415651026022
$ws.Range("Q10").Value = 9.931696522271114
$ws.Range("R10").Value = 89.38526870044002
$ws.Range("S10").Value = 0.0821809639735097
$ws.Range("T10").Value = 0.08218096397350973

# row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "L1cam"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.441487333333333
$ws.Range("H11").Value = 10.324462
$ws.Range("I11").Value = 0.2676541983690312
$ws.Range("J11").Value = 0.2676541983690313
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.165953666666667
$ws.Range("N11").Value = 9.497861
$ws.Range("O11").Value = 0.3368406220840099
$ws.Range("P11").Value = 0.3368406220840099
$ws.Range("Q11").Value = 10.89558944175356
$ws.Range("R11").Value = 98.060304975782
$ws.Range("S11").Value = 0.09015680668202147
$ws.Range("T11").Value = 0.09015680668202149

# row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "L1cam"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.441487333333333
$ws.Range("H12").Value = 10.324462
$ws.Range("I12").Value = 0.2676541983690312
$ws.Range("J12").Value = 0.2676541983690313
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.327024333333334
$ws.Range("N12").Value = 9.981073
$ws.Range("O12").Value = 0.3539776838580724
$ws.Range("P12").Value = 0.3539776838580724
$ws.Range("Q12").Value = 11.44991210085844
$ws.Range("R12").Value = 103.049208907726
$ws.Range("S12").Value = 0.09474361321355873
$ws.Range("T12").Value = 0.09474361321355876

# row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "L1cam"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.441487333333333
$ws.Range("H13").Value = 10.324462
$ws.Range("I13").Value = 0.2676541983690312
$ws.Range("J13").Value = 0.2676541983690313
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.020115
$ws.Range("N13").Value = 0.060345
$ws.Range("O13").Value = 0.002140128955315263
$ws.Range("P13").Value = 0.002140128955315263
$ws.Range("Q13").Value = 0.06922551771
$ws.Range("R13").Value = 0.6230296593900001
$ws.Range("S13").Value = 0.0005728144999412589
$ws.Range("T13").Value = 0.0005728144999412591

# row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "L1cam"
$ws.Range("C14").Value = "Erbb2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.434707333333333
$ws.Range("H14").Value = 7.304122
$ws.Range("I14").Value = 0.1893540717859783
$ws.Range("J14").Value = 0.1893540717859783
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.885873333333334
$ws.Range("N14").Value = 8.657620000000001
$ws.Range("O14").Value = 0.3070415651026022
$ws.Range("P14").Value = 0.3070415651026022
$ws.Range("Q14").Value = 7.026256967737779
$ws.Range("R14").Value = 63.23631270964002
$ws.Range("S14").Value = 0.05813957055971727
$ws.Range("T14").Value = 0.05813957055971728

# row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "L1cam"
$ws.Range("C15").Value = "Erbb2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.434707333333333
$ws.Range("H15").Value = 7.304122
$ws.Range("I15").Value = 0.1893540717859783
$ws.Range("J15").Value = 0.1893540717859783
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.165953666666667
$ws.Range("N15").Value = 9.497861
$ws.Range("O15").Value = 0.3368406220840099
$ws.Range("P15").Value = 0.3368406220840099
$ws.Range("Q15").Value = 7.708170609226889
$ws.Range("R15").Value = 69.37353548304201
$ws.Range("S15").Value = 0.06378214333452921
$ws.Range("T15").Value = 0.06378214333452921

# row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "L1cam"
$ws.Range("C16").Value = "Erbb2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.434707333333333
$ws.Range("H16").Value = 7.304122
$ws.Range("I16").Value = 0.1893540717859783
$ws.Range("J16").Value = 0.1893540717859783
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.327024333333334
$ws.Range("N16").Value = 9.981073
$ws.Range("O16").Value = 0.3539776838580724
$ws.Range("P16").Value = 0.3539776838580724
$ws.Range("Q16").Value = 8.100330542545112
$ws.Range("R16").Value = 72.90297488290601
$ws.Range("S16").Value = 0.06702711575989578
$ws.Range("T16").Value = 0.06702711575989578

# row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "L1cam"
$ws.Range("C17").Value = "Erbb2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.434707333333333
$ws.Range("H17").Value = 7.304122
$ws.Range("I17").Value = 0.1893540717859783
$ws.Range("J17").Value = 0.1893540717859783
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.020115
$ws.Range("N17").Value = 0.060345
$ws.Range("O17").Value = 0.002140128955315263
$ws.Range("P17").Value = 0.002140128955315263
$ws.Range("Q17").Value = 0.04897413801
$ws.Range("R17").Value = 0.4407672420900001
$ws.Range("S17").Value = 0.000405242131836017
$ws.Range("T17").Value = 0.0004052421318360172
